$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '285.02'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '2.44%'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '19'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '28.59'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '4.27%'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '19'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.099'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '5.31%'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '19'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06670'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '4.77%'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '19'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '7.331'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '4.23%'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '19'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.387'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '2.57%'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '19'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.359'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '5.64%'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '19'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9356'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '4.81%'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '19'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1575'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '3.89%'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '19'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06557'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '12.17%'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '19'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07695'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '2.73%'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '19'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.02871'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-1.49%'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '19'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.08965'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.08%'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '19'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001587'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.52%'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '19'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.04470'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '1.37%'
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '19'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0006447'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.71%'
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '19'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.006146'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.48%'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '19'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.476'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.16%'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '19'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.220'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-1.31%'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '19'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.3198'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '0.89%'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '19'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.1305'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-3.32%'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '19'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.048'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '3.90%'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '19'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1522'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '1.15%'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '19'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001179'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.34%'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '19'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004470'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '4.73%'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '19'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001245'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '5.60%'
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '19'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '-2.52%'
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '19'
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '19'
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '19'
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '19'
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '19'
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '19'
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '19'
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '19'
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '19'
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '19'
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '19'
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '19'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04173'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '3.55%'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '19'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006716'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '0.18%'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '19'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-11.81%'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '19'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002012'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-2.24%'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '19'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01211'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '8.31%'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '19'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005666'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '2.25%'
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '19'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '25.93%'
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '19'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-29.60%'
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '19'
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '19'
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '19'
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '19'
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '19'
